$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.517.17'
$ws.Range('E2').Value = '  +2.72%  '

$ws.Range('D3').Value = '2.193.82'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '257.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.30%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '83.00'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +12.20%  '

$ws.Range('E7').Value = '  +1.38%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.593'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.56%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.57'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +11.95%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0917'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.75%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.09'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.81%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.86%  '

$ws.Range('D14').Value = '2.519.87'
$ws.Range('E14').Value = '  +0.70%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.40'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.84%  '

$ws.Range('D16').Value = '2.214.47'
$ws.Range('E16').Value = '  +1.65%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.776'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.56%  '

$ws.Range('D18').Value = '43.466.19'
$ws.Range('E18').Value = '  +2.84%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000102'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.55%  '

$ws.Range('E20').Value = '  +1.14%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '69.34'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.10%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.36'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +12.70%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '230.98'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.11%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.77'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.25%  '

$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.61'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.15%  '

$ws.Range('B27').Value = 'InjectiveProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '40.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +9.03%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.41'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.91%  '

$ws.Range('E29').Value = '  +2.77%  '

$ws.Range('E30').Value = '  +2.79%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '174.12'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.61%  '

$ws.Range('E32').Value = '  +1.75%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0863'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.98%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.33'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.14%  '

$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.122'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.23%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.111'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +4.13%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0361'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.40%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.47'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +7.58%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '12.42'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.19%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.82'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +10.82%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.10'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.30%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '63.13'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +8.42%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.45'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +6.10%  '

$ws.Range('E44').Value = '  +2.70%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '100.19'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0976'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.76%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.21'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.87%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.18'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.81%  '

$ws.Range('E49').Value = '  +2.43%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.442'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.43%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.49'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +17.78%  '
